# Apply "Add data for 2022-12-03" changes to the carjacking-by-neighborhood workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet (tab) and update the header cell that mirrors it.
$ws.Name = "Through 2022-11-25"
$ws.Range("B1").Value2 = "November 2022 (through November 25)"

# 2. Increment existing cell values by 1 (new carjacking recorded against these
#    neighborhood/month combinations).
$cellsToIncrement = @(
    "BE2", "BP2",
    "AI5", "BE5",
    "B6",
    "B7",
    "M9",
    "X11", "AI11",
    "M12",
    "B17", "BP17",
    "M26", "X26",
    "X30",
    "M59",
    "M84",
    "B97", "X97"
)

foreach ($addr in $cellsToIncrement) {
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 + 1
}

# 3. Set brand-new cell values (these were previously empty).
$cellsToSet = @(
    "BE15",
    "X18",
    "AT56",
    "AT77",
    "X82"
)

foreach ($addr in $cellsToSet) {
    $ws.Range($addr).Value2 = 1
}
